# Rewrite the NATMI LR-pair result table to include the FAPs sending/
# target cluster and updated statistics ("Natmi following Dr Hou advice").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ltb"
$ws.Cells.Item(2, 3).Value = "Ltbr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.013823333333333
$ws.Cells.Item(2, 8).Value = 3.04147
$ws.Cells.Item(2, 9).Value = 0.1808135948909178
$ws.Cells.Item(2, 10).Value = 0.1808135948909177
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 9.667057999999999
$ws.Cells.Item(2, 14).Value = 29.001174
$ws.Cells.Item(2, 15).Value = 0.1512832311431697
$ws.Cells.Item(2, 16).Value = 0.1588228070066391
$ws.Cells.Item(2, 17).Value = 9.800688965086666
$ws.Cells.Item(2, 18).Value = 88.20620068577999
$ws.Cells.Item(2, 19).Value = 0.02735406486971016
$ws.Cells.Item(2, 20).Value = 0.02871732268553685

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ltb"
$ws.Cells.Item(3, 3).Value = "Ltbr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.013823333333333
$ws.Cells.Item(3, 8).Value = 3.04147
$ws.Cells.Item(3, 9).Value = 0.1808135948909178
$ws.Cells.Item(3, 10).Value = 0.1808135948909177
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 20.60908733333333
$ws.Cells.Item(3, 14).Value = 61.82726199999999
$ws.Cells.Item(3, 15).Value = 0.3225189424433408
$ws.Cells.Item(3, 16).Value = 0.3385924756140875
$ws.Cells.Item(3, 17).Value = 20.89397361723778
$ws.Cells.Item(3, 18).Value = 188.04576255514
$ws.Cells.Item(3, 19).Value = 0.05831580940359745
$ws.Cells.Item(3, 20).Value = 0.06122212271879856

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ltb"
$ws.Cells.Item(4, 3).Value = "Ltbr"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.013823333333333
$ws.Cells.Item(4, 8).Value = 3.04147
$ws.Cells.Item(4, 9).Value = 0.1808135948909178
$ws.Cells.Item(4, 10).Value = 0.1808135948909177
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 11.51251933333333
$ws.Cells.Item(4, 14).Value = 34.537558
$ws.Cells.Item(4, 15).Value = 0.1801635123472805
$ws.Cells.Item(4, 16).Value = 0.1891424088112641
$ws.Cells.Item(4, 17).Value = 11.67166072558444
$ws.Cells.Item(4, 18).Value = 105.04494653026
$ws.Cells.Item(4, 19).Value = 0.03257601233568603
$ws.Cells.Item(4, 20).Value = 0.03419951888349225

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Ltb"
$ws.Cells.Item(5, 3).Value = "Ltbr"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.013823333333333
$ws.Cells.Item(5, 8).Value = 3.04147
$ws.Cells.Item(5, 9).Value = 0.1808135948909178
$ws.Cells.Item(5, 10).Value = 0.1808135948909177
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 13.01136333333333
$ws.Cells.Item(5, 14).Value = 39.03409
$ws.Cells.Item(5, 15).Value = 0.2036194555411202
$ws.Cells.Item(5, 16).Value = 0.2137673372377884
$ws.Cells.Item(5, 17).Value = 13.19122374581111
$ws.Cells.Item(5, 18).Value = 118.7210137123
$ws.Cells.Item(5, 19).Value = 0.03681716574612134
$ws.Cells.Item(5, 20).Value = 0.03865204071622366

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Ltb"
$ws.Cells.Item(6, 3).Value = "Ltbr"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.013823333333333
$ws.Cells.Item(6, 8).Value = 3.04147
$ws.Cells.Item(6, 9).Value = 0.1808135948909178
$ws.Cells.Item(6, 10).Value = 0.1808135948909177
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 9.1003655
$ws.Cells.Item(6, 14).Value = 18.200731
$ws.Cells.Item(6, 15).Value = 0.1424148585250888
$ws.Cells.Item(6, 16).Value = 0.09967497133022109
$ws.Cells.Item(6, 17).Value = 9.226162885761667
$ws.Cells.Item(6, 18).Value = 55.35697731457
$ws.Cells.Item(6, 19).Value = 0.02575054253580276
$ws.Cells.Item(6, 20).Value = 0.01802258988686643

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ltb"
$ws.Cells.Item(7, 3).Value = "Ltbr"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.003585333333333333
$ws.Cells.Item(7, 8).Value = 0.010756
$ws.Cells.Item(7, 9).Value = 0.0006394378463856988
$ws.Cells.Item(7, 10).Value = 0.0006394378463856988
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 9.667057999999999
$ws.Cells.Item(7, 14).Value = 29.001174
$ws.Cells.Item(7, 15).Value = 0.1512832311431697
$ws.Cells.Item(7, 16).Value = 0.1588228070066391
$ws.Cells.Item(7, 17).Value = 0.03465962528266666
$ws.Cells.Item(7, 18).Value = 0.311936627544
$ws.Cells.Item(7, 19).Value = 0.000096736223516458315872801332
$ws.Cells.Item(7, 20).Value = 0.0001015573136692568

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Ltb"
$ws.Cells.Item(8, 3).Value = "Ltbr"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.003585333333333333
$ws.Cells.Item(8, 8).Value = 0.010756
$ws.Cells.Item(8, 9).Value = 0.0006394378463856988
$ws.Cells.Item(8, 10).Value = 0.0006394378463856988
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 20.60908733333333
$ws.Cells.Item(8, 14).Value = 61.82726199999999
$ws.Cells.Item(8, 15).Value = 0.3225189424433408
$ws.Cells.Item(8, 16).Value = 0.3385924756140875
$ws.Cells.Item(8, 17).Value = 0.07389044778577776
$ws.Cells.Item(8, 18).Value = 0.6650140300719999
$ws.Cells.Item(8, 19).Value = 0.000206230817974563
$ws.Cells.Item(8, 20).Value = 0.0002165088434090743

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Ltb"
$ws.Cells.Item(9, 3).Value = "Ltbr"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.003585333333333333
$ws.Cells.Item(9, 8).Value = 0.010756
$ws.Cells.Item(9, 9).Value = 0.0006394378463856988
$ws.Cells.Item(9, 10).Value = 0.0006394378463856988
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 11.51251933333333
$ws.Cells.Item(9, 14).Value = 34.537558
$ws.Cells.Item(9, 15).Value = 0.1801635123472805
$ws.Cells.Item(9, 16).Value = 0.1891424088112641
$ws.Cells.Item(9, 17).Value = 0.04127621931644444
$ws.Cells.Item(9, 18).Value = 0.371485973848
$ws.Cells.Item(9, 19).Value = 0.0001152033683326283
$ws.Cells.Item(9, 20).Value = 0.0001209448145504781

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Ltb"
$ws.Cells.Item(10, 3).Value = "Ltbr"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.003585333333333333
$ws.Cells.Item(10, 8).Value = 0.010756
$ws.Cells.Item(10, 9).Value = 0.0006394378463856988
$ws.Cells.Item(10, 10).Value = 0.0006394378463856988
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 13.01136333333333
$ws.Cells.Item(10, 14).Value = 39.03409
$ws.Cells.Item(10, 15).Value = 0.2036194555411202
$ws.Cells.Item(10, 16).Value = 0.2137673372377884
$ws.Cells.Item(10, 17).Value = 0.04665007467111111
$ws.Cells.Item(10, 18).Value = 0.41985067204
$ws.Cells.Item(10, 19).Value = 0.0001302019861334424
$ws.Cells.Item(10, 20).Value = 0.0001366909257509368

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Ltb"
$ws.Cells.Item(11, 3).Value = "Ltbr"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.003585333333333333
$ws.Cells.Item(11, 8).Value = 0.010756
$ws.Cells.Item(11, 9).Value = 0.0006394378463856988
$ws.Cells.Item(11, 10).Value = 0.0006394378463856988
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 9.1003655
$ws.Cells.Item(11, 14).Value = 18.200731
$ws.Cells.Item(11, 15).Value = 0.1424148585250888
$ws.Cells.Item(11, 16).Value = 0.09967497133022109
$ws.Cells.Item(11, 17).Value = 0.03262784377266667
$ws.Cells.Item(11, 18).Value = 0.195767062636
$ws.Cells.Item(11, 19).Value = 0.000091065450428606735757720847
$ws.Cells.Item(11, 20).Value = 0.000063735949005952840617887645

# Row 12
$ws.Cells.Item(12, 1).Value = "M1"
$ws.Cells.Item(12, 2).Value = "Ltb"
$ws.Cells.Item(12, 3).Value = "Ltbr"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.664473
$ws.Cells.Item(12, 8).Value = 4.993418999999999
$ws.Cells.Item(12, 9).Value = 0.2968558099164587
$ws.Cells.Item(12, 10).Value = 0.2968558099164587
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 9.667057999999999
$ws.Cells.Item(12, 14).Value = 29.001174
$ws.Cells.Item(12, 15).Value = 0.1512832311431697
$ws.Cells.Item(12, 16).Value = 0.1588228070066391
$ws.Cells.Item(12, 17).Value = 16.090557030434
$ws.Cells.Item(12, 18).Value = 144.815013273906
$ws.Cells.Item(12, 19).Value = 0.04490930610778446
$ws.Cells.Item(12, 20).Value = 0.04714747300716125

# Row 13
$ws.Cells.Item(13, 1).Value = "M1"
$ws.Cells.Item(13, 2).Value = "Ltb"
$ws.Cells.Item(13, 3).Value = "Ltbr"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.664473
$ws.Cells.Item(13, 8).Value = 4.993418999999999
$ws.Cells.Item(13, 9).Value = 0.2968558099164587
$ws.Cells.Item(13, 10).Value = 0.2968558099164587
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 20.60908733333333
$ws.Cells.Item(13, 14).Value = 61.82726199999999
$ws.Cells.Item(13, 15).Value = 0.3225189424433408
$ws.Cells.Item(13, 16).Value = 0.3385924756140875
$ws.Cells.Item(13, 17).Value = 34.30326942097533
$ws.Cells.Item(13, 18).Value = 308.7294247887779
$ws.Cells.Item(13, 19).Value = 0.09574162187241766
$ws.Cells.Item(13, 20).Value = 0.1005131435800387

# Row 14
$ws.Cells.Item(14, 1).Value = "M1"
$ws.Cells.Item(14, 2).Value = "Ltb"
$ws.Cells.Item(14, 3).Value = "Ltbr"
$ws.Cells.Item(14, 4).Value = "M1"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.664473
$ws.Cells.Item(14, 8).Value = 4.993418999999999
$ws.Cells.Item(14, 9).Value = 0.2968558099164587
$ws.Cells.Item(14, 10).Value = 0.2968558099164587
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 11.51251933333333
$ws.Cells.Item(14, 14).Value = 34.537558
$ws.Cells.Item(14, 15).Value = 0.1801635123472805
$ws.Cells.Item(14, 16).Value = 0.1891424088112641
$ws.Cells.Item(14, 17).Value = 19.16227759231133
$ws.Cells.Item(14, 18).Value = 172.460498330802
$ws.Cells.Item(14, 19).Value = 0.05348258537524585
$ws.Cells.Item(14, 20).Value = 0.05614802295721773

# Row 15
$ws.Cells.Item(15, 1).Value = "M1"
$ws.Cells.Item(15, 2).Value = "Ltb"
$ws.Cells.Item(15, 3).Value = "Ltbr"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.664473
$ws.Cells.Item(15, 8).Value = 4.993418999999999
$ws.Cells.Item(15, 9).Value = 0.2968558099164587
$ws.Cells.Item(15, 10).Value = 0.2968558099164587
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 13.01136333333333
$ws.Cells.Item(15, 14).Value = 39.03409
$ws.Cells.Item(15, 15).Value = 0.2036194555411202
$ws.Cells.Item(15, 16).Value = 0.2137673372377884
$ws.Cells.Item(15, 17).Value = 21.65706296152333
$ws.Cells.Item(15, 18).Value = 194.91356665371
$ws.Cells.Item(15, 19).Value = 0.06044561838940757
$ws.Cells.Item(15, 20).Value = 0.06345807602940842

# Row 16
$ws.Cells.Item(16, 1).Value = "M1"
$ws.Cells.Item(16, 2).Value = "Ltb"
$ws.Cells.Item(16, 3).Value = "Ltbr"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.664473
$ws.Cells.Item(16, 8).Value = 4.993418999999999
$ws.Cells.Item(16, 9).Value = 0.2968558099164587
$ws.Cells.Item(16, 10).Value = 0.2968558099164587
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 9.1003655
$ws.Cells.Item(16, 14).Value = 18.200731
$ws.Cells.Item(16, 15).Value = 0.1424148585250888
$ws.Cells.Item(16, 16).Value = 0.09967497133022109
$ws.Cells.Item(16, 17).Value = 15.1473126648815
$ws.Cells.Item(16, 18).Value = 90.883875989289
$ws.Cells.Item(16, 19).Value = 0.0422766781716031
$ws.Cells.Item(16, 20).Value = 0.02958909434263258

# Row 17
$ws.Cells.Item(17, 1).Value = "M2"
$ws.Cells.Item(17, 2).Value = "Ltb"
$ws.Cells.Item(17, 3).Value = "Ltbr"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 2.925126666666667
$ws.Cells.Item(17, 8).Value = 8.775380000000002
$ws.Cells.Item(17, 9).Value = 0.5216911573462379
$ws.Cells.Item(17, 10).Value = 0.5216911573462378
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 9.667057999999999
$ws.Cells.Item(17, 14).Value = 29.001174
$ws.Cells.Item(17, 15).Value = 0.1512832311431697
$ws.Cells.Item(17, 16).Value = 0.1588228070066391
$ws.Cells.Item(17, 17).Value = 28.27736914401334
$ws.Cells.Item(17, 18).Value = 254.49632229612
$ws.Cells.Item(17, 19).Value = 0.07892312394215864
$ws.Cells.Item(17, 20).Value = 0.08285645400027172

# Row 18
$ws.Cells.Item(18, 1).Value = "M2"
$ws.Cells.Item(18, 2).Value = "Ltb"
$ws.Cells.Item(18, 3).Value = "Ltbr"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 2.925126666666667
$ws.Cells.Item(18, 8).Value = 8.775380000000002
$ws.Cells.Item(18, 9).Value = 0.5216911573462379
$ws.Cells.Item(18, 10).Value = 0.5216911573462378
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 20.60908733333333
$ws.Cells.Item(18, 14).Value = 61.82726199999999
$ws.Cells.Item(18, 15).Value = 0.3225189424433408
$ws.Cells.Item(18, 16).Value = 0.3385924756140875
$ws.Cells.Item(18, 17).Value = 60.28419093439556
$ws.Cells.Item(18, 18).Value = 542.55771840956
$ws.Cells.Item(18, 19).Value = 0.1682552803493512
$ws.Cells.Item(18, 20).Value = 0.1766407004718411

# Row 19
$ws.Cells.Item(19, 1).Value = "M2"
$ws.Cells.Item(19, 2).Value = "Ltb"
$ws.Cells.Item(19, 3).Value = "Ltbr"
$ws.Cells.Item(19, 4).Value = "M1"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 2.925126666666667
$ws.Cells.Item(19, 8).Value = 8.775380000000002
$ws.Cells.Item(19, 9).Value = 0.5216911573462379
$ws.Cells.Item(19, 10).Value = 0.5216911573462378
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 11.51251933333333
$ws.Cells.Item(19, 14).Value = 34.537558
$ws.Cells.Item(19, 15).Value = 0.1801635123472805
$ws.Cells.Item(19, 16).Value = 0.1891424088112641
$ws.Cells.Item(19, 17).Value = 33.67557730244889
$ws.Cells.Item(19, 18).Value = 303.0801957220401
$ws.Cells.Item(19, 19).Value = 0.09398971126801599
$ws.Cells.Item(19, 20).Value = 0.09867392215600362

# Row 20
$ws.Cells.Item(20, 1).Value = "M2"
$ws.Cells.Item(20, 2).Value = "Ltb"
$ws.Cells.Item(20, 3).Value = "Ltbr"
$ws.Cells.Item(20, 4).Value = "M2"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 2.925126666666667
$ws.Cells.Item(20, 8).Value = 8.775380000000002
$ws.Cells.Item(20, 9).Value = 0.5216911573462379
$ws.Cells.Item(20, 10).Value = 0.5216911573462378
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 13.01136333333333
$ws.Cells.Item(20, 14).Value = 39.03409
$ws.Cells.Item(20, 15).Value = 0.2036194555411202
$ws.Cells.Item(20, 16).Value = 0.2137673372377884
$ws.Cells.Item(20, 17).Value = 38.05988585602223
$ws.Cells.Item(20, 18).Value = 342.5389727042001
$ws.Cells.Item(20, 19).Value = 0.1062264694194578
$ws.Cells.Item(20, 20).Value = 0.1115205295664054

# Row 21
$ws.Cells.Item(21, 1).Value = "M2"
$ws.Cells.Item(21, 2).Value = "Ltb"
$ws.Cells.Item(21, 3).Value = "Ltbr"
$ws.Cells.Item(21, 4).Value = "sCs"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 2.925126666666667
$ws.Cells.Item(21, 8).Value = 8.775380000000002
$ws.Cells.Item(21, 9).Value = 0.5216911573462379
$ws.Cells.Item(21, 10).Value = 0.5216911573462378
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 9.1003655
$ws.Cells.Item(21, 14).Value = 18.200731
$ws.Cells.Item(21, 15).Value = 0.1424148585250888
$ws.Cells.Item(21, 16).Value = 0.09967497133022109
$ws.Cells.Item(21, 17).Value = 26.61972180046334
$ws.Cells.Item(21, 18).Value = 159.7183308027801
$ws.Cells.Item(21, 19).Value = 0.07429657236725429
$ws.Cells.Item(21, 20).Value = 0.05199955115171612

